$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031443603156359
$ws.Range("D2").Value = 1.039969453658163
$ws.Range("E2").Value = 1.031004155229597
$ws.Range("F2").Value = 1.050738102137029
$ws.Range("I2").Value = 1.037060002326932
$ws.Range("J2").Value = 1.036579299383502
$ws.Range("K2").Value = 1.042753132114885
$ws.Range("L2").Value = 1.033813528208025
$ws.Range("M2").Value = 1.053491546320871
$ws.Range("N2").Value = 1.016125504902895

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032588110826149
$ws.Range("D3").Value = 1.04083911454737
$ws.Range("E3").Value = 1.03198163788482
$ws.Range("F3").Value = 1.051759249168251
$ws.Range("I3").Value = 1.037323290466817
$ws.Range("J3").Value = 1.037364691446546
$ws.Range("K3").Value = 1.043432964176908
$ws.Range("L3").Value = 1.034599030524216
$ws.Range("M3").Value = 1.054324662406226
$ws.Range("N3").Value = 1.016390138025093

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03332857581301
$ws.Range("D4").Value = 1.041401459332506
$ws.Range("E4").Value = 1.032614396645437
$ws.Range("F4").Value = 1.052419888452915
$ws.Range("I4").Value = 1.037491843822578
$ws.Range("J4").Value = 1.037872267439767
$ws.Range("K4").Value = 1.043871840602556
$ws.Range("L4").Value = 1.035106953446458
$ws.Range("M4").Value = 1.054863026032898
$ws.Range("N4").Value = 1.016561042584708

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033639842455116
$ws.Range("D5").Value = 1.041637777062622
$ws.Range("E5").Value = 1.032880471333073
$ws.Range("F5").Value = 1.05269759461534
$ws.Range("I5").Value = 1.03756226984118
$ws.Range("J5").Value = 1.038085503001106
$ws.Range("K5").Value = 1.044056100009341
$ws.Range("L5").Value = 1.035320400754291
$ws.Range("M5").Value = 1.055089182174452
$ws.Range("N5").Value = 1.016632811696215

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033692104067035
$ws.Range("D6").Value = 1.0416774504505
$ws.Range("E6").Value = 1.032925150137057
$ws.Range("F6").Value = 1.052744221142784
$ws.Range("I6").Value = 1.037574069237854
$ws.Range("J6").Value = 1.038121297438294
$ws.Range("K6").Value = 1.044087023653026
$ws.Range("L6").Value = 1.035356234606629
$ws.Range("M6").Value = 1.055127144704022
$ws.Range("N6").Value = 1.016644857404123

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033332735069442
$ws.Range("D7").Value = 1.041404617383286
$ws.Range("E7").Value = 1.032617951701081
$ws.Range("F7").Value = 1.052423599280742
$ws.Range("I7").Value = 1.037492786563199
$ws.Range("J7").Value = 1.037875117290415
$ws.Range("K7").Value = 1.043874303645444
$ws.Range("L7").Value = 1.035109805867638
$ws.Range("M7").Value = 1.054866048617893
$ws.Range("N7").Value = 1.016562001878132

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031830417859163
$ws.Range("D8").Value = 1.040263439074016
$ws.Range("E8").Value = 1.031334445889014
$ws.Range("F8").Value = 1.051083226575734
$ws.Range("I8").Value = 1.03714935685155
$ws.Range("J8").Value = 1.036844855931184
$ws.Range("K8").Value = 1.042983095843988
$ws.Range("L8").Value = 1.034079064906447
$ws.Range("M8").Value = 1.053773250516052
$ws.Range("N8").Value = 1.016215007405816

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029182269696606
$ws.Range("D9").Value = 1.038249599761145
$ws.Range("E9").Value = 1.029074744644662
$ws.Range("F9").Value = 1.048720466072182
$ws.Range("I9").Value = 1.036530318145549
$ws.Range("J9").Value = 1.035024602155403
$ws.Range("K9").Value = 1.041404856476353
$ws.Range("L9").Value = 1.032260080235918
$ws.Range("M9").Value = 1.051842098878015
$ws.Range("N9").Value = 1.01560102435885

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027416183114159
$ws.Range("D10").Value = 1.036905066488834
$ws.Range("E10").Value = 1.027569608114961
$ws.Range("F10").Value = 1.047144713396352
$ws.Range("I10").Value = 1.036108302167639
$ws.Range("J10").Value = 1.033807849717614
$ws.Range("K10").Value = 1.040347434480953
$ws.Range("L10").Value = 1.031045600392612
$ws.Range("M10").Value = 1.050550960147769
$ws.Range("N10").Value = 1.015189995581772

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026651277011971
$ws.Range("D11").Value = 1.036322399753142
$ws.Range("E11").Value = 1.026918179034678
$ws.Range("F11").Value = 1.04646225529266
$ws.Range("I11").Value = 1.035923353211542
$ws.Range("J11").Value = 1.033280205961725
$ws.Range("K11").Value = 1.039888308530279
$ws.Range("L11").Value = 1.030519280193435
$ws.Range("M11").Value = 1.049991001977276
$ws.Range("N11").Value = 1.015011609822848

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026367128463902
$ws.Range("D12").Value = 1.036105899832338
$ws.Range("E12").Value = 1.026676254633032
$ws.Range("F12").Value = 1.046208737570287
$ws.Range("I12").Value = 1.035854322333725
$ws.Range("J12").Value = 1.033084097508294
$ws.Range("K12").Value = 1.039717579805325
$ws.Range("L12").Value = 1.03032371447963
$ws.Range("M12").Value = 1.049782875065798
$ws.Range("N12").Value = 1.014945288052013

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026428080599342
$ws.Range("D13").Value = 1.036152343041364
$ws.Range("E13").Value = 1.026728146192968
$ws.Range("F13").Value = 1.046263119001361
$ws.Range("I13").Value = 1.035869144743628
$ws.Range("J13").Value = 1.033126168789733
$ws.Range("K13").Value = 1.039754210249857
$ws.Range("L13").Value = 1.030365667026981
$ws.Range("M13").Value = 1.049827525041789
$ws.Range("N13").Value = 1.01495951707552

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026627789791531
$ws.Range("D14").Value = 1.036304505262944
$ws.Range("E14").Value = 1.026898180565412
$ws.Range("F14").Value = 1.046441299897521
$ws.Range("I14").Value = 1.035917653883103
$ws.Range("J14").Value = 1.033263997987818
$ws.Range("K14").Value = 1.039874199896064
$ws.Range("L14").Value = 1.030503116039285
$ws.Range("M14").Value = 1.049973800867363
$ws.Range("N14").Value = 1.015006128896965

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026750833413268
$ws.Range("D15").Value = 1.036398247942362
$ws.Range("E15").Value = 1.027002950351564
$ws.Range("F15").Value = 1.046551080050304
$ws.Range("I15").Value = 1.035947497889527
$ws.Range("J15").Value = 1.033348903442531
$ws.Range("K15").Value = 1.039948104436003
$ws.Range("L15").Value = 1.030587794024922
$ws.Range("M15").Value = 1.050063908529631
$ws.Range("N15").Value = 1.015034839844418

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027466943516504
$ws.Range("D16").Value = 1.036943726137846
$ws.Range("E16").Value = 1.027612847704257
$ws.Range("F16").Value = 1.047190002771233
$ws.Range("I16").Value = 1.03612052998293
$ws.Range("J16").Value = 1.033842851171138
$ws.Range("K16").Value = 1.040377878684606
$ws.Range("L16").Value = 1.031080521167991
$ws.Range("M16").Value = 1.050588104021309
$ws.Range("N16").Value = 1.015201825859334

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027916091511164
$ws.Range("D17").Value = 1.037285762904127
$ws.Range("E17").Value = 1.027995501530264
$ws.Range("F17").Value = 1.047590742573713
$ws.Range("I17").Value = 1.036228475834881
$ws.Range("J17").Value = 1.034152481771146
$ws.Range("K17").Value = 1.040647128428707
$ws.Range("L17").Value = 1.031389477096521
$ws.Range("M17").Value = 1.050916680420765
$ws.Range("N17").Value = 1.015306462566801

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028178054770012
$ws.Range("D18").Value = 1.037485221350229
$ws.Range("E18").Value = 1.028218726488903
$ws.Range("F18").Value = 1.047824473180114
$ws.Range("I18").Value = 1.036291225283121
$ws.Range("J18").Value = 1.034333008657705
$ws.Range("K18").Value = 1.040804056157433
$ws.Range("L18").Value = 1.031569643209221
$ws.Range("M18").Value = 1.051108247954979
$ws.Range("N18").Value = 1.015367456070099

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028267374604541
$ws.Range("D19").Value = 1.037553223761045
$ws.Range("E19").Value = 1.028294845441132
$ws.Range("F19").Value = 1.047904166888578
$ws.Range("I19").Value = 1.036312585005646
$ws.Range("J19").Value = 1.034394550888987
$ws.Range("K19").Value = 1.040857543927865
$ws.Range("L19").Value = 1.031631067977839
$ws.Range("M19").Value = 1.051173553037346
$ws.Range("N19").Value = 1.015388246613209

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027867903968174
$ws.Range("D20").Value = 1.037249070335998
$ws.Range("E20").Value = 1.027954443352489
$ws.Range("F20").Value = 1.047547748450366
$ws.Range("I20").Value = 1.036216916354615
$ws.Range("J20").Value = 1.034119269141586
$ws.Range("K20").Value = 1.04061825300654
$ws.Range("L20").Value = 1.031356333464536
$ws.Range("M20").Value = 1.050881436130625
$ws.Range("N20").Value = 1.015295240106021

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026568981206931
$ws.Range("D21").Value = 1.036259699258018
$ws.Range("E21").Value = 1.026848108432137
$ws.Range("F21").Value = 1.046388830685668
$ws.Range("I21").Value = 1.035903378332691
$ws.Range("J21").Value = 1.033223413988367
$ws.Range("K21").Value = 1.03983887115318
$ws.Range("L21").Value = 1.030462642583295
$ws.Range("M21").Value = 1.049930729970796
$ws.Range("N21").Value = 1.01499240457207

$ws.Range("B22").Value = 1.019999999999999
$ws.Range("C22").Value = 1.025752132071722
$ws.Range("D22").Value = 1.035637228314008
$ws.Range("E22").Value = 1.026152773989001
$ws.Range("F22").Value = 1.045660043337849
$ws.Range("I22").Value = 1.03570432016442
$ws.Range("J22").Value = 1.032659470887818
$ws.Range("K22").Value = 1.039347749744975
$ws.Range("L22").Value = 1.029900356015168
$ws.Range("M22").Value = 1.049332210637749
$ws.Range("N22").Value = 1.014801644796063

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026185175373226
$ws.Range("D23").Value = 1.035967251256062
$ws.Range("E23").Value = 1.026521359230097
$ws.Range("F23").Value = 1.046046399715076
$ws.Range("I23").Value = 1.035810027103038
$ws.Range("J23").Value = 1.032958492774691
$ws.Range("K23").Value = 1.039608206217013
$ws.Range("L23").Value = 1.030198471672406
$ws.Range("M23").Value = 1.049649570405062
$ws.Range("N23").Value = 1.014902803859098

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027889677897242
$ws.Range("D24").Value = 1.037265650269487
$ws.Range("E24").Value = 1.027972995684454
$ws.Range("F24").Value = 1.047567175687507
$ws.Range("I24").Value = 1.036222140245873
$ws.Range("J24").Value = 1.034134276731998
$ws.Range("K24").Value = 1.040631300940843
$ws.Range("L24").Value = 1.031371309776747
$ws.Range("M24").Value = 1.050897361772295
$ws.Range("N24").Value = 1.015300311174396

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029866989880697
$ws.Range("D25").Value = 1.038770573432945
$ws.Range("E25").Value = 1.02965869596433
$ws.Range("F25").Value = 1.049331397895954
$ws.Range("I25").Value = 1.036691998024648
$ws.Range("J25").Value = 1.035495752607322
$ws.Range("K25").Value = 1.041813796220098
$ws.Range("L25").Value = 1.032730651530097
$ws.Range("M25").Value = 1.051842098878015
$ws.Range("N25").Value = 1.01576005431269
